# Analysis_parameters.xlsx — add an "Analysis method" column ahead of the
# existing parameters (conditions, filtering_type, mass_spec, comparison,
# control) so users can record whether the data came from DiaNN, MaxQuant or
# Proteome Discoverer.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new column at the front; everything that used to live in
# A:E (conditions/filtering_type/mass_spec/comparison/control) slides right
# to B:F, and their existing data validations slide with them automatically.
$ws.Columns("A:A").Insert()

# New header + sample value for the inserted column.
$ws.Range("A1").Value = "Analysis method"
$ws.Range("A2").Value = "Proteome Discoverer"

# Re-apply the (auto-fit-like) column widths across A:F.
$ws.Columns("A:A").ColumnWidth = 14.92
$ws.Columns("B:B").ColumnWidth = 9.42
$ws.Columns("C:C").ColumnWidth = 12.42
$ws.Columns("D:D").ColumnWidth = 9.6
$ws.Columns("E:E").ColumnWidth = 10.42
$ws.Columns("F:F").ColumnWidth = 6.42

# New dropdown validation for the "Analysis method" cell.
$ws.Range("A2").Validation.Add(3, 1, 1, """DiaNN, MaxQuant, Proteome Discoverer""")

# Match the saved selection/active cell.
[void]$ws.Range("A5").Select()
